# "Corrected small errors in QWET exercise"
# - Fix the k_min value on the Parameters sheet (4.19E-6 -> 4.9E-6).
# - Remove the two erroneous turbulence-parameter rows (cSigTm, cTmOpt),
#   which shifts the remaining cMuMax row up.
# - Leave the Parameters sheet active/selected, matching the author's
#   final view state.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Parameters")

# 1. Correct k_min's value (was 4.1899999999999997E-6).
$ws3.Range("C4").Value = 0.0000049

# 2. Delete rows 7 ("cSigTm", 10, "degree Celsius") and 8 ("cTmOpt", 38,
#    "degree Celsius") entirely; this shifts the old row 9
#    ("cMuMax", "user-defined", "/day") up to become the new row 7.
$ws3.Rows("7:8").Delete()

# 3. Parameters becomes the active/selected sheet.
[void]$ws3.Activate()

# 4. Leave the selection on C5, as in the saved workbook.
[void]$ws3.Range("C5").Select()
